$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix CanClone flag for the "City/SelectScene" row (M6): 0 -> 1
$ws.Range("M6").Value = 1

# Update the active selection saved in the sheet view
$ws.Range("N12").Select() | Out-Null
